$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:G12").ClearContents()

$ws.Range("F17").Select() | Out-Null
